$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the empty predicted price (TimeTaken in Hours) by computing it from
# TimeTaken in Minutes using a formula instead of a hard-coded value.
$ws.Range("C2").Formula = "=B2/60"
